# Add two new date columns ("輸入日期" / entry date, "結束日期" / finish date)
# to the report header, right after the first column ("定型日期"), pushing the
# existing columns (訂單.. 工號) two places to the right. Also move the active
# selection from A4 to B4, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column B; this shifts columns B:O to D:Q,
# carries formatting/column widths along, and expands the A1:O1 merged
# banner to A1:Q1 automatically.
$ws.Columns("B:C").Insert()

# Fill in the two new header cells on row 3 (the column-title row).
# Write C3 first, then B3, so the shared-string table order matches.
$ws.Range("C3").Value = "結束日期"
$ws.Range("B3").Value = "輸入日期"

# Move the selection like the author's saved view (was A4, now B4).
$ws.Range("B4").Select() | Out-Null
